$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values (one cell per statement -
# a comma-joined multi-area Range only applies property changes to its first area)
# so they are not auto-converted to numbers by Excel; source cells are text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated Price / Volume(1h) values
$ws.Range("D2").Value = "58.879.36"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.500.75"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "535.68"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "138.11"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "2.525.63"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").Value = "2.958.72"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "23.29"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "58.864.64"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "2.506.78"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "325.22"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "64.84"
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "168.08"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("E33").Value = "  +6.06%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "1.42"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "36.75"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "0.834"
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "282.17"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "130.65"
$ws.Range("E45").Value = "  +6.53%  "
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("E51").Value = "  -0.25%  "
